$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "CASO01"
$ws.Range("B5").Value = (Get-Date -Year 2020 -Month 6 -Day 27)
$ws.Range("C5").Value = "No permite seleccionar tipo de contrato"

$ws.Range("A6").Value = "CASO02"
$ws.Range("B6").Value = (Get-Date -Year 2020 -Month 6 -Day 27)
$ws.Range("C6").Value = "No aparece perfil de empleado creado "

$ws.Range("A7").Value = "CASO03"
$ws.Range("B7").Value = (Get-Date -Year 2020 -Month 6 -Day 27)
$ws.Range("C7").Value = "Editamos perfil ya creado pero no guarda cambio en el campo cliente"

$ws.Range("A8").Value = "CASO04"
$ws.Range("B8").Value = (Get-Date -Year 2020 -Month 6 -Day 27)
$ws.Range("C8").Value = "Orden Alfabético en el listado de clientes "
